# Upload new version with timestamp
# Inserts two new product rows (DIGENORM SYRUP 120 ML, INDERAL 10 MG 50 TABS)
# into the product table, keeping alphabetical order, and renumbers / retotals
# the remainder of the table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: insert one blank row before row 7 (old "HEALSEC" row), then
#    one more blank row before what is now row 9 (old "PANTOLOC" row, after
#    the first insert shifted it from 8 -> 9).
# ---------------------------------------------------------------------------
$ws.Rows(7).Insert()
$ws.Rows(9).Insert()

# ---------------------------------------------------------------------------
# 2) Copy the formatting (styles only, not values) from neighbouring rows
#    onto the two new blank rows so they keep the same cell styles used by
#    the rest of the table (borders/fill/font/alignment), without minting
#    brand-new style records.
# ---------------------------------------------------------------------------
$ws.Range("A6:N6").Copy()
$ws.Range("A7:N7").PasteSpecial(-4122)

$ws.Range("A8:N8").Copy()
$ws.Range("A9:N9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row heights for the two new rows (matches the other product rows).
$ws.Rows(7).RowHeight = 25.5
$ws.Rows(9).RowHeight = 25.5

# Re-create the merged cell groups for the two new rows.
$ws.Range("B7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()

$ws.Range("B9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()

# ---------------------------------------------------------------------------
# 3) Fill in the two new product rows.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "DIGENORM SYRUP 120 ML"
$ws.Range("H7").Value = "4:0"
$ws.Range("L7").Value = 55
$ws.Range("N7").Value = "1:0"

$ws.Range("B9").Value = "INDERAL 10 MG 50 TABS"
$ws.Range("H9").Value = "0:1"
$ws.Range("L9").Value = 45
$ws.Range("N9").Value = "0:5"

# ---------------------------------------------------------------------------
# 4) Renumber column A (the "م" / index column) sequentially for every
#    product row, now that two rows were inserted in the middle.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 7
$ws.Range("A11").Value = 8
$ws.Range("A12").Value = 9
$ws.Range("A13").Value = 10
$ws.Range("A14").Value = 11
$ws.Range("A15").Value = 12
$ws.Range("A16").Value = 13
$ws.Range("A17").Value = 14
$ws.Range("A18").Value = 15

# ---------------------------------------------------------------------------
# 5) Update the grand total (column K of the totals row) to include the two
#    new rows' price column values (446 + 55 + 45 = 546).
# ---------------------------------------------------------------------------
$ws.Range("K19").Value = 546

Write-Output "edit applied"
